$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D7", "D10", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D41", "D42", "D43", "D44", "D47", "D48")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "94.481.27"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "3.083.29"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "237.19"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "611.63"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "1.11"
$ws.Range("E7").Value = "  +2.21%  "
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "0.814"
$ws.Range("E10").Value = "  +10.75%  "
$ws.Range("D11").Value = "3.078.35"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").Value = "94.007.61"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "0.0000240"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").Value = "34.00"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "3.647.94"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "3.065.76"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "3.58"
$ws.Range("E19").Value = "  -4.68%  "
$ws.Range("D20").Value = "14.43"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "5.78"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").Value = "441.18"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").Value = "8.83"
$ws.Range("E23").Value = "  -5.37%  "
$ws.Range("D24").Value = "0.0000190"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D25").Value = "8.34"
$ws.Range("E25").Value = "  +6.26%  "
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("D27").Value = "84.66"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "11.89"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").Value = "3.236.38"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "0.250"
$ws.Range("E31").Value = "  +10.00%  "
$ws.Range("D32").Value = "0.178"
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("D34").Value = "9.03"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  +30.83%  "
$ws.Range("D36").Value = "7.72"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").Value = "25.42"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").Value = "477.27"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").Value = "0.439"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("D44").Value = "1.27"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -5.93%  "
$ws.Range("D47").Value = "161.45"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "0.673"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("E51").Value = "  +0.15%  "
